# Updated the hazards, added some additional assumptions
#
# The "Assumptions" sheet is removed from the workbook entirely (its
# content - the A-1..A-5 style assumption rows - is gone), leaving only
# the "Hazards" sheet behind. The Hazards sheet's own data is unchanged;
# row heights for two of its rows are tweaked slightly.

$wb = $excel.ActiveWorkbook

# Deleting a worksheet pops a confirmation dialog in interactive Excel;
# suppress it so the automation can proceed unattended.
$excel.DisplayAlerts = $false | Out-Null

$wb.Worksheets("Assumptions").Delete() | Out-Null

$ws = $wb.Worksheets("Hazards")

# Minor row-height adjustments on the surviving Hazards sheet.
$ws.Rows(3).RowHeight = 123.75
$ws.Rows(5).RowHeight = 101.25

# Make sure Hazards (now the only sheet) is the active/selected tab,
# since the sheet that used to hold that status has been removed.
$ws.Activate()

$excel.DisplayAlerts = $true | Out-Null
